$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "5e7fc9b0003c00e10219031e1bb470f100e3021f02f31be7710700e402a3037b197770c400e3032c03bb18a870ba"

$ws.Range("A2:AL2").Select()
